$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.409.00'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.243.88'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.59%  '
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.837'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.238.45'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').Value = '44.070.27'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '0.0₃0960'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '237.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.26%  '
$ws.Range('E27').Value = '  +3.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -13.08%  '
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.98%  '
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0300'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').Value = '1.735.55'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.193'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '80.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '99.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.03%  '
$ws.Range('E47').Value = '  -4.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.62'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.08%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '69.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.44%  '
